# "Generate Report for Handback"
#
# The localization-status report previously showed the zh-cn / de-de rows as
# "In Translation" with no target/handback file recorded. This handback run
# fills in the handback file name + datetime for each localized doc and
# flips the status to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status columns for zh-cn (E) and de-de (F) on both rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Helper: convert a desired OOXML-stored column width (character units) into
# the value to hand to the COM ColumnWidth setter (Excel stores width on a
# Calibri-11 / MDW-7 pixel grid, offset by the 5px padding term).
function Set-GridColumnWidth($ws, [int]$colIndex, [double]$storedWidth) {
    $sixths = [Math]::Round($storedWidth * 6)
    $ws.Columns.Item($colIndex).ColumnWidth = ($sixths - 5) / 6
}

# Overview: widen the now-longer status columns (E, F).
Set-GridColumnWidth $overview 5 29.9777047293527
Set-GridColumnWidth $overview 6 29.9777047293527

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de): fill in Latest Target File (I),
# Latest Handback File (J) and Latest Handback DateTime (K) for both rows,
# update the Status column (C), and hyperlink the new Latest Target File
# cells the same way the Source File Name column already is.
# ---------------------------------------------------------------------------
$mdUrl7ce = "https://github.com/OpenLocalizationTestOrg/oltest/blob/3d11db9f530a2d574b2890789057e8d0250fd83e/e2e/7ce2214b-b03e-4953-8b9e-01e005630c7c.md"
$mdUrlB1d = "https://github.com/OpenLocalizationTestOrg/oltest/blob/3d11db9f530a2d574b2890789057e8d0250fd83e/e2e/b1d92c72-1542-4109-a270-0381627fde3d.md"

function Update-LanguageSheet($ws, [string]$handbackDatetime, [string]$xlfSuffix) {
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl7ce, "", "", "7ce2214b-b03e-4953-8b9e-01e005630c7c.md")
    $ws.Range("J2").Value = "7ce2214b-b03e-4953-8b9e-01e005630c7c.745da7c3acb79213a268b1886cdb9b31c02ba97c.$xlfSuffix"
    $ws.Range("K2").Value = $handbackDatetime

    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrlB1d, "", "", "b1d92c72-1542-4109-a270-0381627fde3d.md")
    $ws.Range("J3").Value = "b1d92c72-1542-4109-a270-0381627fde3d.c777f7efe83c1efadb7c0685f00e790b4c900898.$xlfSuffix"
    $ws.Range("K3").Value = $handbackDatetime

    Set-GridColumnWidth $ws 3 29.9777047293527
    Set-GridColumnWidth $ws 9 40
    Set-GridColumnWidth $ws 10 40
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $zhcn "2016-08-12 14:25:49" "zh-cn.xlf"

$dede = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $dede "2016-08-12 14:25:57" "de-de.xlf"
